$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")
$ws.Range("BN13").Value = "s"
